$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "tba" placeholder for the Vice President of Pledge Education
# with the newly added officer's name.
$ws.Range("A4").Value = "Alex Pham"

# Update the active selection to reflect the cell that was just edited.
$ws.Range("A4").Select()
